# Työaikaraportti (work-time report) update:
#   - A new daily log entry is added just above the existing "Yht"
#     (total) row; the total row moves down one row and its SUM
#     formula is widened to include the new entry.
#   - The active cell selection ends up on H22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Push the current totals row's formatting down to row 24 - it will
#    become the new "Yht" row once we're done.
$ws.Range("B23:D23").Copy() | Out-Null
$ws.Range("B24:D24").PasteSpecial(-4122) | Out-Null

# 2) Give row 23 the same look as the other data rows above it (date /
#    number / wrapped-text styles), copied from the row directly above.
$ws.Range("B22:D22").Copy() | Out-Null
$ws.Range("B23:D23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row heights: the new entry wraps onto two lines like the other
# multi-line entries, and the totals row keeps its usual height.
$ws.Rows.Item(23).RowHeight = 37.5
$ws.Rows.Item(24).RowHeight = 18.75

# New entry: 2024-02-19, 5 hours, description of the work done.
$ws.Range("B23").Value = 45341
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = "Suunitelin ikonin ja kehitin sen. Tein muunnin työkaluun paranuksia."

# Re-establish the totals row on row 24 with a SUM range that now
# covers the newly inserted row.
$ws.Range("B24").Value = "Yht"
$ws.Range("C24").Formula = "=SUM(C6:C23)"

# Move the selection like the author's session ended up.
$ws.Range("H22").Select() | Out-Null
